$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 23 with the new condition entry
$ws.Range("A23").Value = "cityHaveItem"
$ws.Range("B23").Value = "城市有道具卖"
$ws.Range("C23").Value = "city"
$ws.Range("D23").Value = "sellItemNumber"
$ws.Range("E23").Value = ">"
$ws.Range("F23").Value = "number"
$ws.Range("G23").Value = 0

# Update the selected cell in the sheet view
$ws.Range("G20").Select()
